$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateValue = $ws.Range("A13").Value2
$priceValue = $ws.Range("B13").Value2

$ws.Range("A14").Value2 = $dateValue
$ws.Range("B14").Value2 = $priceValue
